$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# --- B2: append " order By ss.study_subject_id ASC LIMIT 100 " (with trailing space) ---
$b2 = $ws.Range("B2").Value2
$b2New = $b2 + "`n order By ss.study_subject_id ASC LIMIT 100 "
$ws.Range("B2").Value2 = $b2New

# --- B3: append " order By samp.sample_id ASC LIMIT 100" ---
$b3 = $ws.Range("B3").Value2
$b3New = $b3 + "`n order By samp.sample_id ASC LIMIT 100"
$ws.Range("B3").Value2 = $b3New

# --- B4: replace trailing "    order by f.file_name" with "     order By f.file_name ASC LIMIT 100" ---
$b4 = $ws.Range("B4").Value2
$oldTail = '    order by f.file_name'
$newTail = '     order By f.file_name ASC LIMIT 100'
$b4New = $b4.Substring(0, $b4.Length - $oldTail.Length) + $newTail
$ws.Range("B4").Value2 = $b4New

# --- Row heights grow because the wrapped text got longer (rows 2 & 3; row 4 was already at Excel's 409.6pt cap) ---
$ws.Rows.Item(2).RowHeight = 360
$ws.Rows.Item(3).RowHeight = 374.4

# --- Selection moved from B4 to C3 ---
$ws.Range("C3").Select() | Out-Null
